# Atualizacao de bases das ligas, do dia: 20-02-2024 as 23:00
# Romania Liga I - fixes to match ids / rows that were in the wrong order,
# plus newly played fixtures appended at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the odds / market data for each fixture row (everything
# except the row-local A/C/D/E/H/I/J columns which are either the running
# index, the (unchanged) league identifiers/date, or the match result which
# stays attached to the row it was originally reported against).
$swapCols = @("B","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Get-RowValues($row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# --- Rows 69 / 70: the two fixtures were swapped (same two teams' data,
#     just reported against the other row's id/result). ---
$row69 = Get-RowValues 69 $swapCols
$row70 = Get-RowValues 70 $swapCols
Set-RowValues 69 $swapCols $row70
Set-RowValues 70 $swapCols $row69

# --- Rows 139 / 140 / 141: three-way rotation (139 <- 140 <- 141 <- 139). ---
# Note rows 139-141 also have their H/I/J (result) columns re-shuffled with
# the rest of the row, unlike the simple two-row swaps above.
$rotCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
$row139 = Get-RowValues 139 $rotCols
$row140 = Get-RowValues 140 $rotCols
$row141 = Get-RowValues 141 $rotCols
Set-RowValues 139 $rotCols $row140
Set-RowValues 140 $rotCols $row141
Set-RowValues 141 $rotCols $row139

# --- Rows 148 / 149: another simple two-row swap. ---
$row148 = Get-RowValues 148 $swapCols
$row149 = Get-RowValues 149 $swapCols
Set-RowValues 148 $swapCols $row149
Set-RowValues 149 $swapCols $row148

# --- New fixtures appended as rows 359-366. ---
$newRows = @(
    @{ A=357; B=6836260; C="Romania Liga I"; D="Romania Liga I"; E=45345.5;    F="AFC Hermannstadt";      G="FC U Craiova 1948";  K=2.1;   L=3.25; M=3.6;   N=2.15;  O=3.3;  P=3.4;   Q=-0.25; R=1.875; S=1.975; T=2.25; U=2.05;  V=1.8;   W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=358; B=6836263; C="Romania Liga I"; D="Romania Liga I"; E=45345.625;  F="CFR Cluj";              G="Dinamo Bucharest";   K=1.5;   L=4.1;  M=6.5;   N=1.444; O=4;    P=8.5;   Q=-1.25; R=2.025; S=1.825; T=2.5;  U=2.05;  V=1.8;   W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=359; B=6836261; C="Romania Liga I"; D="Romania Liga I"; E=45346.375;  F="ACS UTA Batrana Doamna"; G="Petrolul Ploiesti"; K=2.4;   L=3.1;  M=3.1;   N=2.2;   O=3.1;  P=3.4;   Q=-0.25; R=1.9;   S=1.95;  T=2.25; U=2.05;  V=1.8;   W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=360; B=6836265; C="Romania Liga I"; D="Romania Liga I"; E=45346.5;    F="Farul Constanta";       G="Otelul Galati";      K=1.833; L=3.4;  M=4.5;   N=1.8;   O=3.4;  P=4.75;  Q=-0.5;  R=1.8;   S=2.05;  T=2.25; U=2.05;  V=1.8;   W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=361; B=6835220; C="Romania Liga I"; D="Romania Liga I"; E=45346.625;  F="CSM Politehnica Iasi";  G="Rapid Bucuresti";    K=4.2;   L=3.3;  M=1.909; N=4.2;   O=3.3;  P=1.909; Q=0.5;   R=1.9;   S=1.95;  T=2.5;  U=1.975; V=1.875; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=362; B=6836264; C="Romania Liga I"; D="Romania Liga I"; E=45347.375;  F="ACS Sepsi";             G="Universitatea Cluj"; K=2.25;  L=3.2;  M=3.25;  N=2.2;   O=3.2;  P=3.4;   Q=-0.25; R=1.9;   S=1.95;  T=2.25; U=1.85;  V=2;     W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=363; B=6836262; C="Romania Liga I"; D="Romania Liga I"; E=45347.625;  F="FCSB";                  G="FC Botosani";        K=1.363; L=4.75; M=8.5;   N=1.333; O=4.75; P=10;    Q=-1.5;  R=2.05;  S=1.8;   T=2.5;  U=1.975; V=1.875; W=0; X=0; Y=0; Z=0; AA=0 },
    @{ A=364; B=6835795; C="Romania Liga I"; D="Romania Liga I"; E=45348.625;  F="CS U Craiova";          G="FC Voluntari";       K=1.533; L=4.1;  M=6;     N=1.533; O=4.2;  P=6;     Q=-1;    R=1.975; S=1.875; T=2.25; U=1.8;   V=2.05;  W=0; X=0; Y=0; Z=0; AA=0 }
)

$newRowCols = @("A","B","C","D","E","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

$r = 359
foreach ($rowData in $newRows) {
    foreach ($c in $newRowCols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
    $r++
}

# Match the formatting (bold/centered/bordered id column, date-time number
# format on the date column) used by every other data row in the sheet.
$ws.Range("A358").Copy()
$ws.Range("A359:A366").PasteSpecial(-4122)
$ws.Range("E358").Copy()
$ws.Range("E359:E366").PasteSpecial(-4122)
